$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.505.31'
$ws.Range('E2').Value = '  -2.62%  '

$ws.Range('D3').Value = '1.986.67'
$ws.Range('E3').Value = '  -1.50%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -10.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.598'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.35%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.62'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.48%  '

$ws.Range('E9').Value = '  -3.42%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.23%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0747'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.44%  '

$ws.Range('E12').Value = '  -3.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.48%  '

$ws.Range('E14').Value = '  -1.73%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.96'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.73%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.753'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.63%  '

$ws.Range('E17').Value = '  -4.20%  '

$ws.Range('D18').Value = '1.996.92'
$ws.Range('E18').Value = '  -1.07%  '

$ws.Range('D19').Value = '36.430.21'
$ws.Range('E19').Value = '  -2.45%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.90%  '

$ws.Range('D21').Value = '0.0₃0802'
$ws.Range('E21').Value = '  -4.89%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.68%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '221.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.90%  '

$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.37'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.42%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -10.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.65%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.25%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.127'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.22%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '18.86'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.33%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.14%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.117'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.76%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.38'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.71%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0606'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.67%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.25'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.47%  '

$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.79'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.11%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.23'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.21%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.52%  '

$ws.Range('E41').Value = '  -1.07%  '

$ws.Range('D42').Value = '1.454.91'
$ws.Range('E42').Value = '  +3.91%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0923'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.96%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0202'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.67%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -10.78%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.92%  '

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '14.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.46%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.993'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.17%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.59%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.77'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.14%  '

$ws.Range('D51').Value = '2.166.86'
$ws.Range('E51').Value = '  -1.74%  '

